$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3, 55277, 8187, 9468),
    @(4, 29560, 3341, 3474),
    @(5, 95639, 8464, 8386),
    @(6, 2154, 836, 210),
    @(7, 61015, 10107, 8041),
    @(8, 6506, 1712, 1357),
    @(9, 7527, 1409, 872),
    @(10, 3460, 516, 320),
    @(11, 408, 301, 30),
    @(12, 3, 0, 0),
    @(13, 1235, 304, 340),
    @(14, 3856, 1575, 1168),
    @(15, 6302, 2410, 1020),
    @(16, 4510, 1680, 750),
    @(17, 2822, 1110, 241),
    @(18, 21916, 3390, 4136),
    @(19, 1901, 758, 508),
    @(20, 23731, 2999, 3936),
    @(21, 348, 478, 28),
    @(22, 21657, 2953, 3771),
    @(23, 1384, 528, 215),
    @(24, 25403, 3311, 4768),
    @(25, 99052, 9400, 11942),
    @(26, 7637, 2466, 1217),
    @(27, 0, 0, 0),
    @(28, 6768, 1408, 1561),
    @(29, 1776, 536, 376),
    @(30, 18748, 3360, 3302),
    @(31, 595, 211, 277),
    @(32, 3394, 2166, 501),
    @(33, 20794, 4022, 3531),
    @(34, 13114, 3759, 2755),
    @(35, 7249, 800, 1713),
    @(36, 72833, 7096, 7204),
    @(37, 10559, 3441, 1534),
    @(38, 31767, 2472, 3635),
    @(39, 1289, 1184, 218),
    @(40, 2510, 581, 912),
    @(41, 3426, 409, 151),
    @(42, 11803, 650, 366),
    @(43, 329, 123, 72),
    @(44, 993, 72, 89),
    @(45, 1045, 14, 7),
    @(46, 4092, 1123, 521),
    @(47, 15493, 4391, 2706),
    @(48, 39001, 4357, 5367),
    @(49, 18863, 4452, 1622),
    @(50, 14170, 1543, 2112),
    @(51, 40271, 3769, 5800),
    @(52, 6179, 781, 1519),
    @(53, 17292, 3749, 2934),
    @(54, 2517, 1657, 947),
    @(55, 2685, 1573, 182),
    @(56, 4917, 1263, 1474),
    @(57, 15903, 5969, 3278),
    @(58, 17596, 1268, 685),
    @(59, 857775, 129197, 119586)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.Range("B3").Select()
